$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '58.098.29'
$ws.Range("E2").Value = '  +0.01%  '

$ws.Range("D3").Value = '2.354.98'
$ws.Range("E3").Value = '  +0.02%  '

$ws.Range("E4").Value = '  +0.01%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '542.89'
$ws.Range("E5").Value = '  -0.45%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '133.61'
$ws.Range("E6").Value = '  -0.74%  '

$ws.Range("E7").Value = '  +0.05%  '

$ws.Range("E8").Value = '  +4.62%  '

$ws.Range("E9").Value = '  +3.32%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '5.54'
$ws.Range("E10").Value = '  +2.27%  '

$ws.Range("E11").Value = '  -2.18%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.356'
$ws.Range("E12").Value = '  -1.25%  '

$ws.Range("D13").Value = '2.774.95'
$ws.Range("E13").Value = '  +0.81%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '23.79'
$ws.Range("E14").Value = '  +0.82%  '

$ws.Range("D15").Value = '58.049.94'
$ws.Range("E15").Value = '  +0.00%  '

$ws.Range("E16").Value = '  +1.63%  '

$ws.Range("D17").Value = '2.343.99'
$ws.Range("E17").Value = '  -0.59%  '

$ws.Range("E18").Value = '  +2.19%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '4.31'
$ws.Range("E19").Value = '  +2.33%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '329.66'
$ws.Range("E20").Value = '  -1.25%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.82'
$ws.Range("E21").Value = '  +1.11%  '

$ws.Range("E22").Value = '  +0.66%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '63.49'
$ws.Range("E23").Value = '  +2.63%  '

$ws.Range("E24").Value = '  -2.72%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '1.00'
$ws.Range("E25").Value = '  +0.35%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '8.25'
$ws.Range("E26").Value = '  -2.79%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '1.33'
$ws.Range("E27").Value = '  -5.55%  '

$ws.Range("E28").Value = '  -0.27%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '170.52'
$ws.Range("E29").Value = '  +0.41%  '

$ws.Range("E30").Value = '  +0.31%  '

$ws.Range("E31").Value = '  -0.55%  '

$ws.Range("E32").Value = '  -0.84%  '

$ws.Range("B33").Value = 'SuiNetwork'
$ws.Range("C33").Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.01'
$ws.Range("E33").Value = '  -3.11%  '

$ws.Range("B34").Value = 'USDe'
$ws.Range("C34").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.999'
$ws.Range("E34").Value = '  -0.02%  '

$ws.Range("E35").Value = '  -0.02%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '4.18'
$ws.Range("E36").Value = '  -0.75%  '

$ws.Range("E37").Value = '  -2.73%  '

$ws.Range("E38").Value = '  -2.68%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.394'
$ws.Range("E39").Value = '  +3.68%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '141.89'
$ws.Range("E40").Value = '  -5.23%  '

$ws.Range("E41").Value = '  +0.79%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '289.21'
$ws.Range("E42").Value = '  +0.27%  '

$ws.Range("E43").Value = '  +2.21%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0515'
$ws.Range("E44").Value = '  +1.80%  '

$ws.Range("E45").Value = '  +0.94%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '18.93'
$ws.Range("E46").Value = '  -2.10%  '

$ws.Range("E47").Value = '  +2.19%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.383'
$ws.Range("E48").Value = '  +0.62%  '

$ws.Range("E49").Value = '  +0.16%  '

$ws.Range("E50").Value = '  +0.59%  '

$ws.Range("E51").Value = '  +0.60%  '
